# Complete the meeting-notes worksheet with the two remaining agenda rows
# (test-phase division of labor, and final collective review/summary),
# matching the "complete doc by lhw" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: 确定测试阶段分工 / 10分钟 -------------------------------------
$ws.Range("A6").NumberFormat = "m/d/yy h:mm"
$ws.Range("A6").Value = (Get-Date -Year 2015 -Month 3 -Day 17 -Hour 9 -Minute 30 -Second 0)
$ws.Range("B6").Value = "刘瀚文"
$ws.Range("C6").Value = "确定测试阶段分工"
$ws.Range("D6").Value = "10分钟"

# --- Row 7: 集体检查产物，最后整理，阶段总结 / 60分钟 ----------------------
$ws.Range("A7").NumberFormat = "m/d/yy h:mm"
$ws.Range("A7").Value = (Get-Date -Year 2015 -Month 3 -Day 20 -Hour 18 -Minute 30 -Second 0)
$ws.Range("B7").Value = "刘瀚文"
$ws.Range("C7").Value = "集体检查产物，最后整理，阶段总结"
$ws.Range("D7").Value = "60分钟"
$ws.Rows.Item(7).RowHeight = 27

# Move the active selection the way it ended up after entering the data.
$ws.Range("D10").Select()
